# Greenhouse 7 postplanting inventory - update counts and remove Dan's notes block

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update inventory counts (rows 2-5) to reflect the latest post-planting tally ---

# Row 2 (VIBCAS)
$ws.Range("B2").Value = 29
$ws.Range("C2").Value = 17
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 9
$ws.Range("H2").Value = 27
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 0
$ws.Range("Q2").Value = 32
$ws.Range("U2").Value = 0

# Row 3 (ALNINC)
$ws.Range("C3").Value = 4
$ws.Range("H3").Value = 12
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 5
$ws.Range("R3").Value = 0

# Row 4 (LONCAN)
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = 0
$ws.Range("F4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 15
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("V4").Value = 0

# Row 5 (SORAME)
$ws.Range("B5").Value = 21
$ws.Range("C5").Value = 18
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 9
$ws.Range("H5").Value = 17
$ws.Range("L5").Value = 4

# --- Remove "Dan's notes:" block (rows 14-19), the remaining notes below shift up to
#     fill the gap left by Excel's shared-string compaction, but keep their row numbers ---
$ws.Range("A14:Y19").ClearContents()

# --- Update the selected cell to match the saved view state ---
[void]$ws.Range("E30").Select()
